# [Feature]: calculate metrics on whole dataset.
# Update the classification metrics sheet values to reflect metrics
# computed on the whole dataset instead of a subset.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metrics")

# Row 2: Vascular congestion -- support changes
$ws.Range("E2").Value = 73

# Row 3: Interstitial edema -- support changes
$ws.Range("E3").Value = 54

# Row 4: Alveolar edema
$ws.Range("B4").Value = 0.8369704749679076
$ws.Range("D4").Value = 0.9112508735150244
$ws.Range("E4").Value = 652

# Row 5: accuracy
$ws.Range("B5").Value = 0.8369704749679076
$ws.Range("C5").Value = 0.8369704749679076
$ws.Range("D5").Value = 0.8369704749679076
$ws.Range("E5").Value = 0.8369704749679076

# Row 6: macro avg
$ws.Range("B6").Value = 0.2789901583226359
$ws.Range("D6").Value = 0.3037502911716748
$ws.Range("E6").Value = 779

# Row 7: weighted avg
$ws.Range("B7").Value = 0.7005195759680047
$ws.Range("C7").Value = 0.8369704749679076
$ws.Range("D7").Value = 0.7626900764207907
$ws.Range("E7").Value = 779
